# Update the "courses" sheet:
#  - department (col C) renamed from "FACULTY OF TECH SCIENCES" to
#    "Automotive" for the individual-qualification rows, and to
#    "Packages" for the bundled-qualification rows.
#  - promotionValidity (col R) promo text cleared out.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("courses")

# Individual qualifications -> department "Automotive"
$ws.Range("C2:C7").Value = "Automotive"

# Bundled qualification packages -> department "Packages"
$ws.Range("C8:C9").Value = "Packages"

# Promotion no longer valid - clear the promotionValidity column
$ws.Range("R2:R9").ClearContents()

# Reflect the manual edit: scroll so column G is left-most and leave the
# just-edited R column selected, matching how the change was made in Excel.
$ws.Range("R2:R9").Select()
$excel.ActiveWindow.ScrollColumn = 7
